$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data.
# D-column "Price" cells are plain-text values in the source data (some
# look like decimals, e.g. "246.49"); force text via NumberFormat="@" so
# Excel does not auto-convert them to numbers, then ClearFormats() to drop
# the temporary number-format style and keep the cell on the default style,
# matching the original (unstyled) cells.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '37.028.73'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  -0.18%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.056.42'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  +0.24%  '
$ws.Range("E4").Value = '  -0.36%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '246.49'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -1.10%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.658'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -1.87%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '58.66'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -0.50%  '
$ws.Range("E8").Value = '  -0.04%  '
$ws.Range("E9").Value = '  -2.93%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0774'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -2.14%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.111'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +2.28%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '15.49'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -2.97%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.887'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +6.44%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.356.23'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +0.14%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.71'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -0.13%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.006.29'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -2.26%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '18.20'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -2.97%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '36.987.24'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -0.18%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '73.82'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -2.16%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0891'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -1.52%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.45'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +0.15%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '238.22'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +0.04%  '
$ws.Range("E23").Value = '  -0.05%  '
$ws.Range("E24").Value = '  +1.51%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '10.10'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +5.49%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '169.81'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +0.35%  '
$ws.Range("E27").Value = '  -2.77%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '20.12'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -0.06%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.47'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +14.18%  '
$ws.Range("E30").Value = '  -1.77%  '
$ws.Range("E31").Value = '  -1.42%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.67'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +3.39%  '
$ws.Range("E33").Value = '  -1.70%  '
$ws.Range("E34").Value = '  -0.16%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.32'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +4.25%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.83'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +5.76%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0848'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -5.71%  '
$ws.Range("E38").Value = '  -0.44%  '
$ws.Range("E39").Value = '  +2.14%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.05'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -1.96%  '
$ws.Range("E41").Value = '  +0.17%  '
$ws.Range("E42").Value = '  +2.07%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0960'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -10.68%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '97.47'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +0.22%  '
$ws.Range("E45").Value = '  -3.87%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.302.40'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +0.68%  '
$ws.Range("E47").Value = '  -5.21%  '
$ws.Range("E48").Value = '  -0.33%  '
$ws.Range("E49").Value = '  -0.41%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.244.44'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +0.26%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '44.52'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +2.45%  '
